# Updated cryptos list on Tue May 16 12:20:20 UTC 2023 with GitHub Actions
# Refresh price (D) and 1h volume change (E) figures; rows 33/34 also swap
# the HuobiToken / ImmutableX entries (ranking order changed upstream).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.186.47"
$ws.Range("E2").Value = "  -1.81%  "

$ws.Range("D3").Value = "1.824.56"
$ws.Range("E3").Value = "  -1.32%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.93%  "

$ws.Range("D5").Value = "'311.75"
$ws.Range("E5").Value = "  -2.42%  "

$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  -0.90%  "

$ws.Range("D7").Value = "'0.4236"
$ws.Range("E7").Value = "  -1.85%  "

$ws.Range("D8").Value = "'0.3682"
$ws.Range("E8").Value = "  -1.62%  "

$ws.Range("D9").Value = "'0.07240"
$ws.Range("E9").Value = "  -1.57%  "

$ws.Range("D10").Value = "'0.8505"
$ws.Range("E10").Value = "  -3.42%  "

$ws.Range("D11").Value = "'20.95"
$ws.Range("E11").Value = "  -3.08%  "

$ws.Range("D12").Value = "1.836.46"
$ws.Range("E12").Value = "  -0.79%  "

$ws.Range("E13").Value = "  -0.72%  "

$ws.Range("D14").Value = "'0.07092"
$ws.Range("E14").Value = "  -0.75%  "

$ws.Range("E15").Value = "  -2.92%  "

$ws.Range("D16").Value = "'89.65"
$ws.Range("E16").Value = "  +2.05%  "

$ws.Range("E17").Value = "  -0.99%  "

$ws.Range("D18").Value = "'0.000008847"
$ws.Range("E18").Value = "  -1.68%  "

$ws.Range("D19").Value = "'1.004"
$ws.Range("E19").Value = "  -0.84%  "

$ws.Range("D20").Value = "'15.02"
$ws.Range("E20").Value = "  -2.96%  "

$ws.Range("D21").Value = "27.241.25"
$ws.Range("E21").Value = "  -1.63%  "

$ws.Range("D22").Value = "'5.114"
$ws.Range("E22").Value = "  -2.53%  "

$ws.Range("D23").Value = "'10.90"
$ws.Range("E23").Value = "  -2.26%  "

$ws.Range("D24").Value = "2.055.64"
$ws.Range("E24").Value = "  -1.71%  "

$ws.Range("D25").Value = "'1.981"
$ws.Range("E25").Value = "  -1.55%  "

$ws.Range("D26").Value = "'152.03"
$ws.Range("E26").Value = "  -2.37%  "

$ws.Range("D27").Value = "'2.182"
$ws.Range("E27").Value = "  +3.21%  "

$ws.Range("D28").Value = "'18.38"
$ws.Range("E28").Value = "  -1.29%  "

$ws.Range("D29").Value = "'5.226"
$ws.Range("E29").Value = "  -3.21%  "

$ws.Range("D30").Value = "'116.69"
$ws.Range("E30").Value = "  -2.87%  "

$ws.Range("D31").Value = "'0.08834"
$ws.Range("E31").Value = "  -1.12%  "

$ws.Range("D32").Value = "'1.191"
$ws.Range("E32").Value = "  -3.58%  "

$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "'3.010"
$ws.Range("E33").Value = "  +2.87%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7470"
$ws.Range("E34").Value = "  -3.96%  "

$ws.Range("D35").Value = "'4.447"
$ws.Range("E35").Value = "  -2.60%  "

$ws.Range("E36").Value = "  -0.88%  "

$ws.Range("D37").Value = "'1.104"
$ws.Range("E37").Value = "  -2.97%  "

$ws.Range("D38").Value = "'0.01967"
$ws.Range("E38").Value = "  -0.34%  "

$ws.Range("D39").Value = "'0.05237"
$ws.Range("E39").Value = "  -1.81%  "

$ws.Range("D40").Value = "'7.229"
$ws.Range("E40").Value = "  +0.56%  "

$ws.Range("D41").Value = "'2.869"
$ws.Range("E41").Value = "  -0.10%  "

$ws.Range("D42").Value = "'0.1697"
$ws.Range("E42").Value = "  +0.95%  "

$ws.Range("D43").Value = "'0.5032"
$ws.Range("E43").Value = "  -2.43%  "

$ws.Range("D44").Value = "'8.604"
$ws.Range("E44").Value = "  -2.60%  "

$ws.Range("D45").Value = "'10.59"
$ws.Range("E45").Value = "  -0.97%  "

$ws.Range("D46").Value = "'106.51"
$ws.Range("E46").Value = "  -2.59%  "

$ws.Range("D47").Value = "'0.4738"
$ws.Range("E47").Value = "  +0.16%  "

$ws.Range("E48").Value = "  -0.95%  "

$ws.Range("D49").Value = "'0.06384"
$ws.Range("E49").Value = "  -1.87%  "

$ws.Range("D50").Value = "'1.663"
$ws.Range("E50").Value = "  -2.22%  "

$ws.Range("D51").Value = "'1.868"
$ws.Range("E51").Value = "  -0.26%  "
